$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.198.43"
$ws.Range("E2").Value = "  +0.91%  "

$ws.Range("D3").Value = "3.089.27"
$ws.Range("E3").Value = "  +0.61%  "

$ws.Range("D5").Value = "'559.99"
$ws.Range("E5").Value = "  +1.72%  "

$ws.Range("D6").Value = "'144.00"
$ws.Range("E6").Value = "  +2.44%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "3.086.94"
$ws.Range("E8").Value = "  +0.71%  "

$ws.Range("D9").Value = "'0.505"
$ws.Range("E9").Value = "  +0.58%  "

$ws.Range("D10").Value = "'0.153"
$ws.Range("E10").Value = "  +1.10%  "

$ws.Range("D11").Value = "'6.12"
$ws.Range("E11").Value = "  -6.55%  "

$ws.Range("D12").Value = "'0.471"
$ws.Range("E12").Value = "  +3.12%  "

$ws.Range("D13").Value = "'0.0000228"
$ws.Range("E13").Value = "  -0.21%  "

$ws.Range("D14").Value = "'35.12"
$ws.Range("E14").Value = "  +0.21%  "

$ws.Range("D15").Value = "3.584.76"
$ws.Range("E15").Value = "  +0.56%  "

$ws.Range("D16").Value = "64.193.79"
$ws.Range("E16").Value = "  +0.88%  "

$ws.Range("D17").Value = "3.078.83"
$ws.Range("E17").Value = "  +0.36%  "

$ws.Range("E18").Value = "  +1.21%  "

$ws.Range("D19").Value = "'6.75"
$ws.Range("E19").Value = "  -0.56%  "

$ws.Range("D20").Value = "'484.21"
$ws.Range("E20").Value = "  -0.22%  "

$ws.Range("D21").Value = "'13.96"
$ws.Range("E21").Value = "  +1.42%  "

$ws.Range("D22").Value = "'0.675"
$ws.Range("E22").Value = "  -0.51%  "

$ws.Range("D23").Value = "'7.56"
$ws.Range("E23").Value = "  +3.66%  "

$ws.Range("D24").Value = "'14.20"
$ws.Range("E24").Value = "  +11.96%  "

$ws.Range("D25").Value = "'81.32"
$ws.Range("E25").Value = "  +0.31%  "

$ws.Range("E26").Value = "  +0.07%  "

$ws.Range("D27").Value = "'2.80"
$ws.Range("E27").Value = "  +1.02%  "

$ws.Range("D28").Value = "'8.01"
$ws.Range("E28").Value = "  +0.67%  "

$ws.Range("D29").Value = "'2.07"
$ws.Range("E29").Value = "  +2.83%  "

$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  -0.07%  "

$ws.Range("D31").Value = "'26.41"
$ws.Range("E31").Value = "  +0.53%  "

$ws.Range("D32").Value = "'1.15"
$ws.Range("E32").Value = "  -0.92%  "

$ws.Range("D33").Value = "'2.49"
$ws.Range("E33").Value = "  +1.07%  "

$ws.Range("D34").Value = "'5.59"
$ws.Range("E34").Value = "  -2.58%  "

$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "'6.22"
$ws.Range("E35").Value = "  +3.36%  "

$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "'55.86"
$ws.Range("E36").Value = "  +0.27%  "

$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").Value = "'2.99"
$ws.Range("E37").Value = "  +15.76%  "

$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").Value = "'450.94"
$ws.Range("E38").Value = "  -3.83%  "

$ws.Range("D39").Value = "'0.0408"
$ws.Range("E39").Value = "  +2.31%  "

$ws.Range("D40").Value = "'0.0819"
$ws.Range("E40").Value = "  -0.96%  "

$ws.Range("D41").Value = "2.972.80"
$ws.Range("E41").Value = "  -2.94%  "

$ws.Range("D42").Value = "'8.25"
$ws.Range("E42").Value = "  -0.38%  "

$ws.Range("E43").Value = "  -5.59%  "

$ws.Range("D44").Value = "'27.96"
$ws.Range("E44").Value = "  -1.33%  "

$ws.Range("D45").Value = "'0.261"
$ws.Range("E45").Value = "  +1.45%  "

$ws.Range("D47").Value = "'2.13"
$ws.Range("E47").Value = "  +2.29%  "

$ws.Range("E48").Value = "  +1.51%  "

$ws.Range("D49").Value = "'119.33"
$ws.Range("E49").Value = "  +1.68%  "

$ws.Range("D50").Value = "0.0₃0515"
$ws.Range("E50").Value = "  -0.27%  "

$ws.Range("D51").Value = "'2.08"
$ws.Range("E51").Value = "  -0.09%  "
